$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the data; existing values in A1:A500 shift
# down to A2:A501 and a new value appears at A501 (continuation of the
# shifted series). Put the "Rewards" header label in the now-empty A1.
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "Rewards"

# Restore/update the view: drop the old frozen/scrolled "B1" top-left cell
# and move the active selection to B14.
$ws.Range("B14").Select()

# Update the line chart: give the series an explicit name (pulled from the
# new header cell) and repoint its value range to the shifted data block.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Rewards_DDQN!`$A`$1,,Rewards_DDQN!`$A`$2:`$A`$501,1)"

# Show a legend under the chart.
$chart.HasLegend = $true
$chart.Legend.Position = -4107

# The chart was anchored starting at row 1; since a row was inserted above
# it, shift the chart down by one row (row height is 14.4 points) so it
# keeps floating over the same data/plot region (now rows 2-33 instead of
# 1-32).
$co.Top = $co.Top + 14.4

Write-Host "done"
